# Microgrid_cmd.xlsx — startup/comms tweak pass.
#
# 1) Move the cursor on the WTGmode sheet (the tab that was left selected)
#    from F12 to A25 — this is the cell the author was last looking at
#    while wiring up the hardcoded startup commands.
# 2) Force a full recalculation so the volatile RAND()-driven V-ramp table
#    on the "V" sheet picks up fresh sampled values, matching a normal
#    save-after-edit recalc pass.

$wb = $excel.ActiveWorkbook

# --- WTGmode: update the saved selection -------------------------------
$wtg = $wb.Worksheets.Item("WTGmode")
$wtg.Activate()
$wtg.Range("A25").Select()

# --- Recalculate the whole workbook -------------------------------------
# "V" sheet (B2:B402) holds volatile `RAND()` formulas (575+0.5*RAND(),
# B53-RAND(), B114+0.1*RAND(), ...); recalculating re-samples them, which
# is the only change on that sheet.
$excel.Calculate()
